$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-14 Sunday" "2024-01-15 Monday"
Replace-Text "71×67=" "76×42="
Replace-Text "72×12=" "48×42="
Replace-Text "88×61=" "84×12="
Replace-Text "81×49=" "84×15="
Replace-Text "80×75=" "30×41="
Replace-Text "62×15=" "43×80="
Replace-Text "19×67=" "53×38="
Replace-Text "77×81=" "82×75="
Replace-Text "99×44=" "25×61="
Replace-Text "87×56=" "87×66="
Replace-Text "60×93=" "97×15="
Replace-Text "24×21=" "67×35="
Replace-Text "51×46=" "85×62="
Replace-Text "79×71=" "27×70="
Replace-Text "21×19=" "47×43="
Replace-Text "55×27=" "96×92="
Replace-Text "53×35=" "59×47="
Replace-Text "38×48=" "71×77="
Replace-Text "99×66=" "84×23="
Replace-Text "69×90=" "59×73="
Replace-Text "97×74=" "14×39="
Replace-Text "91×76=" "47×53="
Replace-Text "80×39=" "46×25="
Replace-Text "67×11=" "21×11="
Replace-Text "40×16=" "51×34="
